$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells keep their original text formatting (avoid Excel
# auto-converting numeric-looking strings like "253.10" into numbers,
# which would silently drop significant trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.638.45"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.098.02"
$ws.Range("E3").Value = "  +9.63%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.10"
$ws.Range("E5").Value = "  +1.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.658"
$ws.Range("E6").Value = "  -6.32%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.49"
$ws.Range("E8").Value = "  +5.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.27"
$ws.Range("E9").Value = "  +3.71%  "

$ws.Range("E10").Value = "  +0.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0747"
$ws.Range("E11").Value = "  -1.95%  "

$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.80"
$ws.Range("E13").Value = "  +0.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.405.38"
$ws.Range("E14").Value = "  +9.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.837"
$ws.Range("E15").Value = "  +2.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.100.17"
$ws.Range("E16").Value = "  +10.02%  "

$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.597.11"
$ws.Range("E18").Value = "  -1.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.18"
$ws.Range("E19").Value = "  -2.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0837"
$ws.Range("E20").Value = "  -2.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.33"
$ws.Range("E21").Value = "  -2.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "240.42"
$ws.Range("E22").Value = "  -4.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.29"
$ws.Range("E23").Value = "  +1.91%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("E25").Value = "  -2.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.45"
$ws.Range("E26").Value = "  +2.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.26"
$ws.Range("E28").Value = "  +4.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.00"
$ws.Range("E29").Value = "  -9.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.37"
$ws.Range("E30").Value = "  +50.10%  "

$ws.Range("E31").Value = "  -4.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.51"
$ws.Range("E32").Value = "  -2.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0620"
$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("E34").Value = "  +20.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.981"
$ws.Range("E35").Value = "  +11.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0898"
$ws.Range("E36").Value = "  -0.72%  "

$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("E38").Value = "  -2.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.10"
$ws.Range("E39").Value = "  -5.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.36"
$ws.Range("E40").Value = "  -10.33%  "

$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("E42").Value = "  +6.43%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.18"
$ws.Range("E43").Value = "  -7.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.56"
$ws.Range("E44").Value = "  -7.62%  "

$ws.Range("E45").Value = "  -3.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.341.25"
$ws.Range("E46").Value = "  -0.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0845"
$ws.Range("E47").Value = "  +3.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.12"
$ws.Range("E48").Value = "  +9.80%  "

$ws.Range("E49").Value = "  +2.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.298.77"
$ws.Range("E50").Value = "  +10.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.26"
$ws.Range("E51").Value = "  -5.95%  "
